$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1818181818181818
$ws.Range("C2").Value = 0.6363636363636364
$ws.Range("P2").Value = 0.1818181818181818
$ws.Range("P3").Value = 0.5714285714285714
$ws.Range("S3").Value = 0.4285714285714285
$ws.Range("J6").Value = 0.2727272727272727
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.6363636363636364
$ws.Range("J7").Value = 0.07692307692307693
$ws.Range("R7").Value = 0.1538461538461539
$ws.Range("S7").Value = 0.7692307692307693
$ws.Range("B8").Value = 0.0131578947368421
$ws.Range("F8").Value = 0.03947368421052631
$ws.Range("J8").Value = 0.131578947368421
$ws.Range("O8").Value = 0.0131578947368421
$ws.Range("Q8").Value = 0.1052631578947368
$ws.Range("R8").Value = 0.09210526315789473
$ws.Range("S8").Value = 0.6052631578947368
$ws.Range("B9").Value = 0.2307692307692308
$ws.Range("F9").Value = 0.07692307692307693
$ws.Range("J9").Value = 0.07692307692307693
$ws.Range("Q9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.5384615384615384
$ws.Range("B10").Value = 0.07462686567164178
$ws.Range("D10").Value = 0.01492537313432836
$ws.Range("F10").Value = 0.05970149253731343
$ws.Range("J10").Value = 0.05970149253731343
$ws.Range("O10").Value = 0.01492537313432836
$ws.Range("Q10").Value = 0.08955223880597014
$ws.Range("R10").Value = 0.08955223880597014
$ws.Range("S10").Value = 0.5970149253731343
$ws.Range("G11").Value = 0.1739130434782609
$ws.Range("K11").Value = 0.2173913043478261
$ws.Range("L11").Value = 0.5652173913043478
$ws.Range("S11").Value = 0.04347826086956522
$ws.Range("G12").Value = 0.6153846153846154
$ws.Range("J12").Value = 0.1538461538461539
$ws.Range("S12").Value = 0.2307692307692308
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5
$ws.Range("H15").Value = 0.2222222222222222
$ws.Range("I15").Value = 0.2222222222222222
$ws.Range("K15").Value = 0.2222222222222222
$ws.Range("O15").Value = 0.1111111111111111
$ws.Range("S15").Value = 0.2222222222222222
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("I16").Value = 0.1428571428571428
$ws.Range("J16").Value = 0.2857142857142857
$ws.Range("S16").Value = 0.2857142857142857
$ws.Range("H17").Value = 0.4666666666666667
$ws.Range("I17").Value = 0.06666666666666667
$ws.Range("J17").Value = 0.1333333333333333
$ws.Range("K17").Value = 0.06666666666666667
$ws.Range("S17").Value = 0.2666666666666667
$ws.Range("H18").Value = 0.4
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.1333333333333333
$ws.Range("K18").Value = 0.1333333333333333
$ws.Range("M18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.2
$ws.Range("F19").Value = 0.01282051282051282
$ws.Range("H19").Value = 0.3846153846153846
$ws.Range("I19").Value = 0.05128205128205128
$ws.Range("J19").Value = 0.25
$ws.Range("K19").Value = 0.08333333333333333
$ws.Range("M19").Value = 0.00641025641025641
$ws.Range("O19").Value = 0.01923076923076923
$ws.Range("S19").Value = 0.1923076923076923
